$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$ws.Range("E5").Value = 161
$ws.Range("F5").Value = 110
$ws.Range("H5").Value = 121
$ws.Range("E6").Value = 51
$ws.Range("F6").Value = 37
$ws.Range("H6").Value = 47
$ws.Range("E7").Value = 43
$ws.Range("F7").Value = 29
$ws.Range("H7").Value = 33
$ws.Range("E9").Value = 13
$ws.Range("E10").Value = 688
$ws.Range("F10").Value = 383
$ws.Range("H10").Value = 478
$ws.Range("E11").Value = 457
$ws.Range("F11").Value = 258
$ws.Range("H11").Value = 323
$ws.Range("E12").Value = 695
$ws.Range("F12").Value = 417
$ws.Range("H12").Value = 503
$ws.Range("E13").Value = 165
$ws.Range("E15").Value = 201
$ws.Range("F15").Value = 95
$ws.Range("H15").Value = 146
$ws.Range("E16").Value = 234
$ws.Range("E18").Value = 62
$ws.Range("F18").Value = 34
$ws.Range("H18").Value = 51
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 45
$ws.Range("H20").Value = 82
$ws.Range("E21").Value = 152
$ws.Range("F21").Value = 90
$ws.Range("H21").Value = 121
$ws.Range("E23").Value = 226
$ws.Range("E24").Value = 264
$ws.Range("E25").Value = 328
$ws.Range("F25").Value = 181
$ws.Range("H25").Value = 241
$ws.Range("E26").Value = 195
$ws.Range("F26").Value = 122
$ws.Range("H26").Value = 147
$ws.Range("E27").Value = 381
$ws.Range("F27").Value = 209
$ws.Range("H27").Value = 291
$ws.Range("E28").Value = 227
$ws.Range("F28").Value = 112
$ws.Range("H28").Value = 164
$ws.Range("E30").Value = 252
$ws.Range("F30").Value = 157
$ws.Range("H30").Value = 209
$ws.Range("E31").Value = 82
$ws.Range("E32").Value = 214
$ws.Range("F32").Value = 137
$ws.Range("H32").Value = 175
$ws.Range("E33").Value = 327
$ws.Range("E34").Value = 250
$ws.Range("F34").Value = 176
$ws.Range("H34").Value = 214
$ws.Range("E35").Value = 184
$ws.Range("F37").Value = 110
$ws.Range("H37").Value = 147
$ws.Range("E39").Value = 198
$ws.Range("F39").Value = 103
$ws.Range("H39").Value = 154
$ws.Range("E40").Value = 305
$ws.Range("E41").Value = 440
$ws.Range("F41").Value = 221
$ws.Range("H41").Value = 313
$ws.Range("E42").Value = 460
$ws.Range("F42").Value = 263
$ws.Range("H42").Value = 324
$ws.Range("E43").Value = 142
$ws.Range("F43").Value = 79
$ws.Range("H43").Value = 106
$ws.Range("E44").Value = 365
$ws.Range("F44").Value = 192
$ws.Range("H44").Value = 260
$ws.Range("E45").Value = 179
$ws.Range("F45").Value = 100
$ws.Range("H45").Value = 139
$ws.Range("E46").Value = 388
$ws.Range("F46").Value = 226
$ws.Range("H46").Value = 290
$ws.Range("E47").Value = 535
$ws.Range("F47").Value = 303
$ws.Range("H47").Value = 395
$ws.Range("E48").Value = 267
$ws.Range("F48").Value = 130
$ws.Range("H48").Value = 174
$ws.Range("E49").Value = 338
$ws.Range("E50").Value = 281
$ws.Range("F50").Value = 155
$ws.Range("H50").Value = 228
$ws.Range("E51").Value = 266
$ws.Range("F51").Value = 132
$ws.Range("H51").Value = 206
$ws.Range("E52").Value = 33